$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the engine name typo/correction: "Gnuchess 5.60" -> "Gnu Chesss 5.60"
$ws.Range("A4").Value = "Gnu Chesss 5.60"

# Update benchmark numbers in the reference table (rows 4-8, cols B/C)
$ws.Range("B4").Value = 2811
$ws.Range("C4").Value = 2817
$ws.Range("C5").Value = 2701
$ws.Range("C6").Value = 2907
$ws.Range("B7").Value = 2648
$ws.Range("C7").Value = 2614
$ws.Range("B8").Value = 2573
$ws.Range("C8").Value = 2518

# Update the Floyd version comparison table (rows 12-15)
$ws.Range("C13").Value = 2295
$ws.Range("B14").Value = 2404
$ws.Range("B15").Value = 2579
$ws.Range("C15").Value = 2652

# Column A width adjustment (engine reports col width 5/6 char wider than the
# ColumnWidth value supplied, so back the target "14.5" out accordingly)
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666

# Sheet view selection
$ws.Range("C9").Select()

# Workbook window geometry (best effort - host may not persist these to the
# saved <workbookView> element, but set them in case it does)
$win = $wb.Windows.Item(1)
$win.Left = 1420
$win.Top = 120
$win.Width = 13560
$win.Height = 16700
